$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("G19").Value = 3.75
$ws.Range("I19").Value = 2.38
$ws.Range("J19").Value = 4.75
$ws.Range("K19").Value = 1.73
$ws.Range("L19").Value = 3.4
$ws.Range("M19").Value = 1.17
$ws.Range("N19").Value = 5
$ws.Range("S19").Value = 1.8
$ws.Range("T19").Value = 2
$ws.Range("U19").Value = 2.63
$ws.Range("V19").Value = 1.44
$ws.Range("W19").Value = 7
$ws.Range("X19").Value = 17
$ws.Range("Y19").Value = 17
$ws.Range("AB19").Value = 67
$ws.Range("AD19").Value = 6
$ws.Range("AH19").Value = 5
$ws.Range("AI19").Value = 9.5
$ws.Range("AK19").Value = 23
$ws.Range("AN19").Value = 5
$ws.Range("AO19").Value = 26
$ws.Range("AQ19").Value = 101
$ws.Range("AT19").Value = 1.91
$ws.Range("AW19").Value = 4
$ws.Range("AX19").Value = 17
$ws.Range("AZ19").Value = 51

# Row 20
$ws.Range("H20").Value = 2.9
$ws.Range("K20").Value = 1.8
$ws.Range("L20").Value = 3.25
$ws.Range("N20").Value = 5
$ws.Range("AC20").Value = 5
$ws.Range("AR20").Value = 151
$ws.Range("AT20").Value = 2

# Row 21
$ws.Range("O21").Value = 1.62
$ws.Range("P21").Value = 2.2

# Row 87
$ws.Range("G87").Value = 2.47
$ws.Range("H87").Value = 3.5
$ws.Range("I87").Value = 2.37
$ws.Range("J87").Value = 3
$ws.Range("L87").Value = 2.87
$ws.Range("W87").Value = 9.75
$ws.Range("X87").Value = 12.5
$ws.Range("Z87").Value = 23
$ws.Range("AA87").Value = 15
$ws.Range("AD87").Value = 6.4
$ws.Range("AH87").Value = 9.75
$ws.Range("AI87").Value = 12.5
$ws.Range("AJ87").Value = 8
$ws.Range("AK87").Value = 22
$ws.Range("AN87").Value = 4.75
$ws.Range("AW87").Value = 4.65
$ws.Range("AY87").Value = 16.5
$ws.Range("BB87").Value = 150

# Row 105
$ws.Range("G105").Value = 2.55
$ws.Range("H105").Value = 3.3
$ws.Range("I105").Value = 2.7
$ws.Range("M105").Value = 1.06
$ws.Range("N105").Value = 10
$ws.Range("O105").Value = 1.3
$ws.Range("P105").Value = 3.4
$ws.Range("Q105").Value = 2.03
$ws.Range("R105").Value = 1.83
$ws.Range("S105").Value = 1.36
$ws.Range("T105").Value = 3
$ws.Range("X105").Value = 13
$ws.Range("Y105").Value = 10
$ws.Range("Z105").Value = 23
$ws.Range("AA105").Value = 21
$ws.Range("AB105").Value = 29
$ws.Range("AC105").Value = 10
$ws.Range("AH105").Value = 9
$ws.Range("AI105").Value = 13
$ws.Range("AJ105").Value = 10
$ws.Range("AK105").Value = 26
$ws.Range("AL105").Value = 21
$ws.Range("AM105").Value = 29
$ws.Range("AN105").Value = 4.5
$ws.Range("AT105").Value = 3
$ws.Range("AU105").Value = 7.5
$ws.Range("AY105").Value = 23

# Row 106
$ws.Range("G106").Value = 11
$ws.Range("I106").Value = 1.33
$ws.Range("M106").Value = 1.08
$ws.Range("N106").Value = 8
$ws.Range("S106").Value = 1.44
$ws.Range("T106").Value = 2.63
$ws.Range("U106").Value = 2.63
$ws.Range("V106").Value = 1.44
$ws.Range("W106").Value = 17
$ws.Range("AD106").Value = 9
$ws.Range("AE106").Value = 29
$ws.Range("AI106").Value = 5.5
$ws.Range("AJ106").Value = 9.5
$ws.Range("AT106").Value = 2.63
$ws.Range("AU106").Value = 11
$ws.Range("AW106").Value = 3.1
